$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161 (this shifts the former rows 161..261
# down to 162..262, growing the used range to A1:R262).
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row with its data.
$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").Value = 44603
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = 100114014
$ws.Range("G161").Value = "Betarraga"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 3000
$ws.Range("K161").Value = 700
$ws.Range("L161").Value = 700
$ws.Range("M161").Value = 700
$ws.Range("N161").Value = "$/paquete 5 unidades"
$ws.Range("O161").Value = "Región del Maule"
$ws.Range("P161").Value = 140
$ws.Range("Q161").Value = 5
$ws.Range("R161").Value = "Hortaliza"
